# Update "想去人数" (want-to-go count) figures across the workbook's sheets,
# matching the values output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3530
$ws1.Range("F4").Value = 144
$ws1.Range("F5").Value = 7020
$ws1.Range("F6").Value = 3203
$ws1.Range("F7").Value = 56
$ws1.Range("F8").Value = 135
$ws1.Range("F13").Value = 15
$ws1.Range("F15").Value = 589
$ws1.Range("F16").Value = 29

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 30

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3530
$ws4.Range("F3").Value = 30
$ws4.Range("F5").Value = 144
$ws4.Range("F6").Value = 7020
$ws4.Range("F7").Value = 3203
$ws4.Range("F8").Value = 56
$ws4.Range("F9").Value = 135
$ws4.Range("F14").Value = 15
$ws4.Range("F16").Value = 589
$ws4.Range("F17").Value = 29
